# "fixed error in data"
# The Week-Start date in row 158 (column C) had been entered with the wrong
# decade (a 2010 date instead of the correct 2020 date). C159:C191 are
# shared formulas ("=C158+7" chained down), so correcting the single
# anchor cell recalculates the whole cascade of week-start dates below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# C158 was 40448 (2010-09-27); the correct value is 44101 (2020-09-27).
$ws.Range("C158").Value = 44101

# Reflect where the author was looking after making/verifying the fix.
$ws.Range("C159").Select()
$excel.ActiveWindow.ScrollRow = 173
$excel.ActiveWindow.ScrollColumn = 1
